$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2) appended below the existing header row.
# Force Text format first so the numeric-looking values ("2300.00", etc.)
# are stored as literal text (shared strings) instead of being
# auto-coerced into numbers, matching the source data.
$ws.Range("A2:H2").NumberFormat = "@"

$ws.Range("A2").Value = "ab445aa0-8f60-4ed3-b2c3-120b767cbc3e"
$ws.Range("B2").Value = "443c417b-e01e-404c-a964-27f3671840fa"
$ws.Range("C2").Value = "ogya test"
$ws.Range("D2").Value = "SAMSUNG, XIAOMI"
$ws.Range("E2").Value = "2300.00"
$ws.Range("F2").Value = "700.00"
$ws.Range("G2").Value = "3000.00"
$ws.Range("H2").Value = "2300.00"

# Drop the Text number-format override again so the row ends up on the
# workbook's default (General) cell style, same as the source data.
$ws.Range("A2:H2").ClearFormats()
